$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-19 18:35:28"

for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
